$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, E, J, M, N, R, S, T, V, X share the same new value across all
# data rows (2-13); column A (Result ID) gets a distinct value per row.
$ws.Range("B2:B13").Value = "1021009"
$ws.Range("E2:E13").Value = "20221111-Cocci-125718"
$ws.Range("J2:J13").Value = "A00375"
$ws.Range("M2:M13").Value = "TestComplexSite_20221110"
$ws.Range("N2:N13").Value = "TestFarm1_20221110"
$ws.Range("R2:R13").Value = "11/11/2022"
$ws.Range("S2:S13").Value = "2:57 PM"
$ws.Range("T2:T13").Value = "CartridgeCocci5718"
$ws.Range("V2:V13").Value = "AFrancisco"
$ws.Range("X2:X13").Value = "11/13/2022"

$resultIds = @(
    "A1109301",
    "A1109302",
    "A1109303",
    "A1109304",
    "A1109305",
    "A1109306",
    "A1109307",
    "A1109308",
    "A1109309",
    "A1109310",
    "A1109311",
    "A1109312"
)

for ($i = 0; $i -lt $resultIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $resultIds[$i]
}
